# DD_KARMEN_FRANZI: removed whitespace and index
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Variables")
$ws.Activate()

# --- "removed index": clear the running index numbers in column A (rows 2-51) ---
$ws.Range("A2:A51").ClearContents()

# --- "removed whitespace": trim the trailing (non-breaking) space off two labels ---
$ws.Range("B14").Value = "Flours_milled_prod_NCI"
$ws.Range("B17").Value = "Bread_NCI"

# --- restore the user's on-screen selection at save time ---
$ws.Range("C34").Select()

$wb.Save()
